{"js": "// Replace each two-digit-divided-by-one-digit expression with its new value.\n// Every \"old\" expression text in this worksheet is unique, so an exact\n// search-and-replace per pair is unambiguous and safe.\nconst replacements = [\n  [\"18\u00f77=\", \"59\u00f75=\"],\n  [\"13\u00f74=\", \"88\u00f72=\"],\n  [\"43\u00f78=\", \"52\u00f75=\"],\n  [\"42\u00f73=\", \"46\u00f77=\"],\n  [\"53\u00f73=\", \"29\u00f78=\"],\n  [\"85\u00f74=\", \"44\u00f79=\"],\n  [\"32\u00f76=\", \"99\u00f79=\"],\n  [\"36\u00f72=\", \"56\u00f78=\"],\n  [\"43\u00f72=\", \"53\u00f77=\"],\n  [\"75\u00f75=\", \"19\u00f76=\"],\n  [\"27\u00f77=\", \"96\u00f76=\"],\n  [\"76\u00f76=\", \"55\u00f75=\"],\n  [\"22\u00f76=\", \"82\u00f72=\"],\n  [\"12\u00f75=\", \"71\u00f73=\"],\n  [\"23\u00f72=\", \"13\u00f77=\"],\n  [\"95\u00f76=\", \"62\u00f74=\"],\n  [\"40\u00f72=\", \"36\u00f74=\"],\n  [\"31\u00f78=\", \"28\u00f73=\"],\n  [\"70\u00f72=\", \"96\u00f75=\"],\n  [\"21\u00f72=\", \"58\u00f74=\"],\n  [\"80\u00f78=\", \"48\u00f74=\"],\n  [\"43\u00f79=\", \"89\u00f77=\"],\n  [\"59\u00f77=\", \"54\u00f76=\"],\n  [\"10\u00f74=\", \"29\u00f72=\"],\n  [\"89\u00f79=\", \"33\u00f77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit-divided-by-one-digit expression with its new value.\n# Every \"old\" expression text in this worksheet is unique, so an exact\n# find & replace-all per pair is unambiguous and safe.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"18\u00f77=\", \"59\u00f75=\"),\n    @(\"13\u00f74=\", \"88\u00f72=\"),\n    @(\"43\u00f78=\", \"52\u00f75=\"),\n    @(\"42\u00f73=\", \"46\u00f77=\"),\n    @(\"53\u00f73=\", \"29\u00f78=\"),\n    @(\"85\u00f74=\", \"44\u00f79=\"),\n    @(\"32\u00f76=\", \"99\u00f79=\"),\n    @(\"36\u00f72=\", \"56\u00f78=\"),\n    @(\"43\u00f72=\", \"53\u00f77=\"),\n    @(\"75\u00f75=\", \"19\u00f76=\"),\n    @(\"27\u00f77=\", \"96\u00f76=\"),\n    @(\"76\u00f76=\", \"55\u00f75=\"),\n    @(\"22\u00f76=\", \"82\u00f72=\"),\n    @(\"12\u00f75=\", \"71\u00f73=\"),\n    @(\"23\u00f72=\", \"13\u00f77=\"),\n    @(\"95\u00f76=\", \"62\u00f74=\"),\n    @(\"40\u00f72=\", \"36\u00f74=\"),\n    @(\"31\u00f78=\", \"28\u00f73=\"),\n    @(\"70\u00f72=\", \"96\u00f75=\"),\n    @(\"21\u00f72=\", \"58\u00f74=\"),\n    @(\"80\u00f78=\", \"48\u00f74=\"),\n    @(\"43\u00f79=\", \"89\u00f77=\"),\n    @(\"59\u00f77=\", \"54\u00f76=\"),\n    @(\"10\u00f74=\", \"29\u00f72=\"),\n    @(\"89\u00f79=\", \"33\u00f77=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
